$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Footnote Test. The footnote anchor should have an inset of 12pt, the content should have an inset of 24pt.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Footnote Test with a link.", 2)
